# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 45188 (2023-09-19) to 45189 (2023-09-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows are 2 through 484 (row 1 is the header row).
$ws.Range("C2:C484").Value = 45189
